# Updates patient data on the "HOJA DE INGRESO Y EGRESO" sheet
# (new patient: FRANCISCO ANTONIO CASTELLANOS LOPEZ) as part of the
# change to add the barcode ("para agregar codigo de barras").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Patient name / record number -------------------------------------
$ws.Range("A6").Value = "CASTELLANOS "
$ws.Range("C6").Value = "LÓPEZ "
$ws.Range("E6").Value = "FRANCISCO "
$ws.Range("G6").Value = "ANTONIO "
$ws.Range("I6").Value = "2017-0034279/201773430"

# --- Current address -----------------------------------------------------
$ws.Range("A8").Value = "4 AV. A 29-39 RESD. VOULEVART SUR"
$ws.Range("D8").Value = "CIUDAD SAN CRISTOBAL "
$ws.Range("F8").Value = "MIXCO "
$ws.Range("H8").Value = "GUATEMALA"
$ws.Range("J8").Value = "2443-6964"

# --- Birth date / age / place / sex (these look like a date & number to
#     Excel's auto-detection, so force them to text and then restore the
#     original cell formatting/type via a formats-only paste from an
#     untouched donor cell with the same original style). ----------------
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1938-08-11"
$ws.Range("B6").Copy()
$ws.Range("A12").PasteSpecial(-4122)

$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "79"
$ws.Range("G12").Copy()
$ws.Range("F12").PasteSpecial(-4122)

$ws.Range("H12").Value = "BELICE"
$ws.Range("J12").Value = "Masculino"

# --- Occupation / cedula --------------------------------------------------
$ws.Range("D14").Value = "NO"

$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "2337250631703"
$ws.Range("D6").Copy()
$ws.Range("H14").PasteSpecial(-4122)

# --- Spouse / parents ------------------------------------------------------
$ws.Range("A16").Value = "ANA MARIA SANTA CRUZ "
$ws.Range("A18").Value = "FRANCISCO PEDRO CASTELLANOS "
$ws.Range("F18").Value = "RUFINA LOPEZ "

# --- Emergency contact -------------------------------------------------
$ws.Range("A20").Value = "ALVARO CASTELLANOS "
$ws.Range("F20").Value = "HIJO "
$ws.Range("H20").Value = ""
$ws.Range("J20").Value = "2443-6964"

# --- Hospitalization dates / service ------------------------------------
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").Value = "13:13:17"
$ws.Range("D24").Value = "CL. 27"

$excel.CutCopyMode = 0
